# Slide 14 ("ORM - Entity Framework") - Content Placeholder bullet text tweak.
#   Paragraph 1: "...to Add/Modify/Delete/Update objects." -> "...to Add/Read/Delete/Update objects."
#   Paragraph 2 (level-1 bullet): "No need to deal with writing SQL" -> "...writing much SQL"
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(14)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# --- Paragraph 1 ---
# Original: "Uses simple common syntax (LINQ) to Add/Modify/Delete/Update objects. "
# Re-key the "Add/Modify/Delete/Update " chunk as "Add/Read/Delete/Update ", which
# splits the single run into three runs (matching the authored edit).
$para1 = $tr.Paragraphs(1, 1)
$oldChunk1 = "Add/Modify/Delete/Update "
$start1 = $para1.Text.IndexOf($oldChunk1) + 1
$midRun1 = $para1.Characters($start1, $oldChunk1.Length)
$midRun1.Text = "Add/Read/Delete/Update "

# --- Paragraph 2 ---
# Original: "No need to deal with writing SQL"
# Insert "much " between "writing " and "SQL" by re-keying the "writing " chunk as
# "writing much ", which splits the single run into three runs.
$para2 = $tr.Paragraphs(2, 1)
$oldChunk2 = "writing "
$start2 = $para2.Text.IndexOf($oldChunk2) + 1
$midRun2 = $para2.Characters($start2, $oldChunk2.Length)
$midRun2.Text = "writing much "
